# Slide 4 title: "But…we aren't measuring O2 or DIC"
# Split the "2" in "O2" into its own run and format it as subscript
# (PowerPoint renders this as baseline="-25000" in the run properties),
# matching the target edit:
#   <a:t>But…we aren't measuring O</a:t>  (unchanged run)
#   <a:t>2</a:t>                          (new run, Font.Subscript = True)
#   <a:t> or DIC</a:t>                    (new run, same formatting as before)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# "But…we aren't measuring O2 or DIC"
#  1234567890123456789012345 6  7890123
#                            O  2
# The "2" is the 26th character (1-based) of the title text.
$two = $tr.Characters(26, 1)
$two.Font.Subscript = $true
